# Update the "Fonds de solidarite" Volet 1 regional/categorie-juridique
# worksheet with the 2022-06-14 data refresh: bump the count of aides
# (column C) and the total amount (column E) for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 11;  C = 1267;    E = 46826598 },
    @{ Row = 13;  C = 187871;  E = 1168646488 },
    @{ Row = 36;  C = 211209;  E = 404259591 },
    @{ Row = 54;  C = 75193;   E = 361072600 },
    @{ Row = 91;  C = 18888;   E = 75401391 },
    @{ Row = 121; C = 1306474; E = 2275750876 },
    @{ Row = 127; C = 9165;    E = 110799628 },
    @{ Row = 129; C = 633944;  E = 3437195477 },
    @{ Row = 132; C = 586127;  E = 3475312863 },
    @{ Row = 136; C = 26710;   E = 144432396 },
    @{ Row = 156; C = 12417;   E = 40823580 },
    @{ Row = 204; C = 265656;  E = 1271520756 },
    @{ Row = 224; C = 39687;   E = 260941034 },
    @{ Row = 240; C = 205944;  E = 1070126718 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$wb.Save()
